$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = -0.077473474454563984
$ws.Range("B1").Value = 0.077347282100561188
$ws.Range("A2").Value = -0.031235099902573893
$ws.Range("B2").Value = 0.03085838360140869
$ws.Range("A3").Value = 0.081157732328303922
$ws.Range("B3").Value = -0.081398187754935236
$ws.Range("A4").Value = -0.20259949976977154
$ws.Range("B4").Value = 0.20154911804906916
$ws.Range("A5").Value = -0.19554911819677123
$ws.Range("B5").Value = 0.19343050589907129
$ws.Range("A6").Value = -0.11150425371007433
$ws.Range("B6").Value = 0.11133175161021791
$ws.Range("A7").Value = -0.09133175178959263
$ws.Range("B7").Value = 0.090898668968018015
$ws.Range("A8").Value = -0.07089866914952303
$ws.Range("B8").Value = 0.070526739018530904
$ws.Range("A9").Value = -0.06452673917692664
$ws.Range("B9").Value = 0.064211777905224032
$ws.Range("A10").Value = -0.056035714789018698
$ws.Range("B10").Value = 0.05598534244858655
$ws.Range("A11").Value = -0.05148534260686688
$ws.Range("B11").Value = 0.051403915744369755
$ws.Range("A12").Value = -0.045403915905932735
$ws.Range("B12").Value = 0.04515529958885578
$ws.Range("A13").Value = -0.039155299753012685
$ws.Range("B13").Value = 0.039087420801925177
$ws.Range("A14").Value = -0.027087420977773391
$ws.Range("B14").Value = 0.027054346507985549
$ws.Range("A15").Value = -0.021054346673412994
$ws.Range("B15").Value = 0.021028287909079957
$ws.Range("A16").Value = -0.01502828807501766
$ws.Range("B16").Value = 0.015004798664170593
$ws.Range("A17").Value = -0.0090047988307775384
$ws.Range("B17").Value = 0.0089999998277496829
$ws.Range("A18").Value = -0.036112182350763078
$ws.Range("B18").Value = 0.036097010346004765
$ws.Range("A19").Value = -0.027097010497110219
$ws.Range("B19").Value = 0.027014136860057825
$ws.Range("A20").Value = -0.018014137012423603
$ws.Range("B20").Value = 0.018004316873055259
$ws.Range("A21").Value = -0.0090043170255880156
$ws.Range("B21").Value = 0.0089999998473393461
$ws.Range("A22").Value = -0.093933939621917872
$ws.Range("B22").Value = 0.093625014765351011
$ws.Range("A23").Value = -0.084625014918114694
$ws.Range("B23").Value = 0.084125049656342377
$ws.Range("A24").Value = -0.042125049871888365
$ws.Range("B24").Value = 0.041999999783368658
$ws.Range("A25").Value = -0.075926252342075173
$ws.Range("B25").Value = 0.075851203569222037
$ws.Range("A26").Value = -0.069851203722723909
$ws.Range("B26").Value = 0.069761675117614885
$ws.Range("A27").Value = -0.063761675271649665
$ws.Range("B27").Value = 0.063480645596647456
$ws.Range("A28").Value = -0.057480645752994164
$ws.Range("B28").Value = 0.057305495039933874
$ws.Range("A29").Value = -0.068749346371243192
$ws.Range("B29").Value = 0.068588859321820905
$ws.Range("A30").Value = -0.048588859506342619
$ws.Range("B30").Value = 0.048257760264846183
$ws.Range("A31").Value = -0.027019806631740195
$ws.Range("B31").Value = 0.027001013461454448
$ws.Range("A32").Value = -0.0060010136497430366
$ws.Range("B32").Value = 0.0059999998390924958
